# Apply edit: add new task "Install Vive at home" to C8, matching style of
# neighboring "Good" cells (e.g. C7), and update the active selection to C9.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set new value in C8
$ws.Range("C8").Value = "Install Vive at home"

# Apply the "Good" cell style (green fill) to match the neighboring task
# cells in column C (e.g. C3, C4, C7 all use this style).
$ws.Range("C8").Style = "Good"

# Update selection to C9 as in the diff
$ws.Range("C9").Select()
